# Auto-generated edit script: refresh crypto price snapshot
# (GitHub Actions symbol-list update, 2022-12-22 18:48 UTC).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds prices stored as literal text (not numbers) in the
# source sheet. Mark each touched cell as Text first so Excel's
# numeric auto-detection doesn't strip significant trailing zeros
# or switch tiny values to scientific notation.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

$ws.Range("D2").Value = "241.89"
$ws.Range("D3").Value = "21.81"
$ws.Range("D4").Value = "5.397"
$ws.Range("D5").Value = "0.05680"
$ws.Range("D6").Value = "3.403"
$ws.Range("D7").Value = "6.293"
$ws.Range("D8").Value = "0.8079"
$ws.Range("D9").Value = "0.9114"
$ws.Range("D10").Value = "0.1428"
$ws.Range("D11").Value = "0.07279"
$ws.Range("D12").Value = "0.03038"
$ws.Range("D13").Value = "0.03120"
$ws.Range("D14").Value = "0.09339"
$ws.Range("D15").Value = "3.907"
$ws.Range("D16").Value = "0.001584"
$ws.Range("D17").Value = "0.04818"
$ws.Range("D18").Value = "0.0005810"
$ws.Range("D19").Value = "0.006308"
$ws.Range("D20").Value = "0.004062"
$ws.Range("D21").Value = "0.0009935"
$ws.Range("D22").Value = "0.0001501"
$ws.Range("D23").Value = "3.736"
$ws.Range("D26").Value = "0.1309"
$ws.Range("D27").Value = "0.0003998"
$ws.Range("D40").Value = "0.03803"
$ws.Range("D41").Value = "0.006682"
$ws.Range("D43").Value = "0.003201"
$ws.Range("D44").Value = "0.006818"
$ws.Range("D45").Value = "0.00005612"
$ws.Range("D46").Value = "0.00000000750"
$ws.Range("D47").Value = "0.5798"
$ws.Range("D49").Value = "0.00002101"
$ws.Range("D50").Value = "0.01010"

# Rows 20/21: BitKan and HotbitToken swapped ranking positions.
$ws.Range("B20").Value = "HotbitToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
$ws.Range("E20").Value = "19HotbitTokenHTB"
$ws.Range("B21").Value = "BitKan"
$ws.Range("C21").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
$ws.Range("E21").Value = "20BitKanKAN"

